$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = "bahasa inggris"
$ws.Range("B2").Value = "big"
$ws.Range("A3").Value = "pendidikan kewarganegaraan"
$ws.Range("B3").Value = "pkn"
$ws.Range("A4").Value = "pendidikan agama islam"
$ws.Range("B4").Value = "pai"
$ws.Range("A5:B5").ClearContents()
$ws.Range("A5").Select() | Out-Null
